$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceOne = 1

# --- 1. Ativação date update -------------------------------------------------
$d.Content.Find.Execute(
    "Ativação: 01/01/2012", $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false, "Ativação: 01/01/2022", $wdReplaceOne)

# --- 2. Objetivos paragraph --------------------------------------------------
$d.Content.Find.Execute(
    "Fornecer uma visão atual dos processos industriais que utilizam a conversão química ou bioquímica como rota de transformação da matéria prima em produto.",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Conferir aos alunos uma noção prática das aplicações, à escala industrial, de processos e produtos estudados em disciplinas de química orgânica e ainda uma visão global das matérias primas mais importantes na área da indústria química de base orgânica.",
    $wdReplaceOne)

$objPara = $d.Paragraphs(6)
$objPara.Range.InsertParagraphAfter()
$objNew = $d.Paragraphs(7).Range
$objNew.MoveEnd(1, -1)
$objNew.Text = "Check the students a practical notion of applications on an industrial scale, processes and products studied in organic chemistry disciplines and also an overview of the most important raw materials in the chemical industry of organic base"
$objNew.Font.Italic = $true

# --- 3. Programa resumido paragraph ------------------------------------------
$d.Content.Find.Execute(
    "Indústria de Papel e Celulose Indústria Sucroalcooleira; Indústrias Biotecnológicas; Tecnologia de Óleos e Gorduras; Indústria de Sabões e Detergentes; Visitas Técnicas.",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Petróleo, Gás Natural e Petroquímica; Química Fina; Processos Unitários Orgânicos: Nitração; Esterificação; Alquilação e Acilação; Hidrogenação; Sulfonação/Sulfatação; Oxidação.",
    $wdReplaceOne)

$resPara = $d.Paragraphs(11)
$resPara.Range.InsertParagraphAfter()
$resNew = $d.Paragraphs(12).Range
$resNew.MoveEnd(1, -1)
$resNew.Text = "Oil, Natural Gas and Petrochemicals; Fine Chemistry; Organic Unit Processes: Nitration; Esterification; Alkylation and Acylation; Hydrogenation; Sulphonation/Sulfation; Oxidation."
$resNew.Font.Italic = $true

# --- 4. Programa (full) paragraph --------------------------------------------
$d.Content.Find.Execute(
    "Indústrias de papel e celulose. Matérias-primas lignocelulósicas. Métodos de separação dos componentes das matérias-primas lignocelulósicas (físico, químico e biotecnológico). Aproveitamento das diversas frações (celulose, hemicelulose, lignina). Descrição do setor industrial de papel e celulose. Processo Kraft. Branqueamento da celulose. Manufatura de papel. Tendências do setor. Indústria açucareira. Terminologia. Produção e consumo do açúcar. Composição da cana de açúcar. Fases do processo de transformação da matéria-prima em açúcar não refinado. Refinação do açúcar. Subprodutos. Aproveitamento global da cana de açúcar. Indústria Alcooleira. Matérias primas para produção de álcool industrial. Fluxograma de obtenção do álcool. Detalhamento das etapas de fermentação e recuperação do produto. Processos industriais de obtenção de álcool. Diferenças entre usinas autônomas e anexa. Aproveitamento dos sub-produtos Indústrias Biotecnológicas. Definição. Terminologia. Áreas de pesquisas em biotecnologia industrial. Classificação dos produtos. Etapas de desenvolvimento de um processo biotecnológico.Exemplos de processos (detalhamento). Tecnologia de óleos e gorduras. Definição, denominação, classificação. Emprego e economia. Propriedades físicas. Extração de óleos vegetais (algodão e soja). Refino de óleos brutos. Aproveitamento de sub-produtos. Modificação de óleos e gorduras por enzimas. Manufatura de sabões e detergentes. Conceitos sobre tensoativos. Tensoativos sintéticos (classificação). Formulação de detergentes. Emprego de enzimas na formulação de detergentes. Fabricação de sabões e sabonetes. Visita Técnica.",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Petróleo, Gás Natural e Petroquímica; 2- Química Fina: Características, Química Fina X Química de Base, Principais Segmentos (Defensivos Agrícolas, Fármacos, Catalisadores, Corantes e Pigmentos, Especialidades); 3- Processos Unitários Orgânicos: 3.1- Nitração; 3.2- Esterificação; 3.3- Alquilação e Acilação; 3.4- Hidrogenação; 3.5- Sulfonação/Sulfatação; 3.6- Oxidação.",
    $wdReplaceOne)

$progPara = $d.Paragraphs(14)
$progPara.Range.InsertParagraphAfter()
$progNew = $d.Paragraphs(15).Range
$progNew.MoveEnd(1, -1)
$progNew.Text = "1- Oil, Natural Gas and Petrochemicals 2- Fine Chemistry: Characteristics, Fine Chemistry vs. Basic Chemistry, Main Segments (Pesticides, Drugs, Catalysts, Dyes and Pigments, Specialties); 3- Organic Unit Processes: 3.1- Nitration; 3.2- Esterification; 3.3- Alkylation and Acylation; 3.4- Hydrogenation; 3.5- Sulphonation/Sulfation; 3.6- Oxidation."
$progNew.Font.Italic = $true

# --- 5. Avaliação: Método ----------------------------------------------------
$d.Content.Find.Execute(
    "Serão oferecidas duas provas (0,0-10,0 pontos). A média final será a média aritmética das notas",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos.",
    $wdReplaceOne)

# --- 6. Avaliação: Critério --------------------------------------------------
$d.Content.Find.Execute(
    "Média final maior que 5,0 - aprovado; - Média final menor que 5,0 e maior que 3,0 - recuperação; - Média final menor que 3,0 - reprovado;",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Provas em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula.",
    $wdReplaceOne)

# --- 7. Avaliação: Norma de recuperação --------------------------------------
$d.Content.Find.Execute(
    "Prova escrita (0,0 - 10,0) para alunos que tenham média final maior ou igual a 3,0 e inferior a 5,0. A nota final será a média aritmética entre a média final e a prova escrita.",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação.",
    $wdReplaceOne)

# --- 8. Bibliografia ----------------------------------------------------------
$d.Content.Find.Execute(
    "1. SHREVE, R.N. & BRINK, J.A. Indústrias de Processos Químicos, Editora Guanabara Dois, 1980.2. CONSIDINE, D. M. (editor-in-chief). Chemical and Process Technology Encyclopedia McGraw-Hill, 1974.3. D ALMEIDA, M. L. O. Celulose e Papel., Vol 1 e 2 (Tecnologia de fabricação da pasta celulósica e Tecnologia de fabricação de papel), Editado pelo ITP/SP, 1989.4. BAYRNA, C. Tecnologia do Açúcar (II). Coleção Canavieira Nº 15.5.FERNANDES, H. Açúcar e Álcool, Coleção Canavieira  Nº 4.6. Manual Econômico da Indústria Química, Editado pelo CEPED/ SPCT/ Governo do Estado da Bahia, 1985.7.URGEL DE ALMEIDA LIMA, U. A.; AQUARONE, E.; BORZANI, W. Biotecnologia: Tecnologia das Fermentações, Editora Edgard Blucher Ltda, 1975.",
    $false, $false, $false, $false, $false, $true,
    $wdFindContinue, $false,
    "Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.Manual econômico da indústria química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.T.W. Graham Solomons, Craig B. Fryhle Hoboken, NJ. Organic chemistry; John Wiley, 9th ed; c2008.Revistas:Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;Química & Derivados, São Paulo, SP: QD, v. 1, n. 1, nov. 1965-; Disponível em: http://www.quimica.com.br/pquimica/category/revista/BiodieselBR. Curitiba: BiodieselBR, v.1, n.1, out/dez.2007 -;Petróleo & Energia, São Paulo, SP: , v. 1, n. 1, ; Disponível em: http://www.petroleoenergia.com.br/petroleo/category/revista-petroleo-e-energia/.",
    $wdReplaceOne)

Write-Output "done"
